# Elimna EC anteriores y se agregan nuevos, se modifica base de datos
#
# The "Estado de Cuenta" worker/period detail table (rows 16-19, cols B:G)
# is reorganised so each worker's two mora-period rows are adjacent:
#   CESAR ANDRES DUEÑAS D ELIA -> period 1806 then 1805
#   ANGEL ALEXANDER PELOCHE TANG -> period 1806 then 1805

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16: CESAR / CC / 1235045954, period 1806 (was 1805)
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "1235045954"
$ws.Range("D16").Value = "CESAR ANDRES DUEÑAS D ELIA"
$ws.Range("E16").Value = "1806"
$ws.Range("F16").Value = 44000
$ws.Range("G16").Value = 1100000

# Row 17: CESAR / CC / 1235045954, period 1805 (was ANGEL / PE / 1805)
$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "1235045954"
$ws.Range("D17").Value = "CESAR ANDRES DUEÑAS D ELIA"
$ws.Range("E17").Value = "1805"
$ws.Range("F17").Value = 44000
$ws.Range("G17").Value = 1100000

# Row 18: ANGEL / PE / 927514822121998, period 1806 (was ANGEL / PE / 1806, unchanged person, swapped order)
$ws.Range("B18").Value = "PE"
$ws.Range("C18").Value = "927514822121998"
$ws.Range("D18").Value = "ANGEL ALEXANDER PELOCHE TANG"
$ws.Range("E18").Value = "1806"
$ws.Range("F18").Value = 40000
$ws.Range("G18").Value = 1000000

# Row 19: ANGEL / PE / 927514822121998, period 1805 (was CESAR / CC / 1806)
$ws.Range("B19").Value = "PE"
$ws.Range("C19").Value = "927514822121998"
$ws.Range("D19").Value = "ANGEL ALEXANDER PELOCHE TANG"
$ws.Range("E19").Value = "1805"
$ws.Range("F19").Value = 40000
$ws.Range("G19").Value = 1000000
